$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(82.3517499923706, 75.71989986419678, 42.23194306946534, 4.88114871383457, 0.003224937955573992),
    @(82.47575016021729, 75.82839992523193, 36.63498098307812, 4.257566297211689, -0.01797219810061368),
    @(82.66000022888184, 75.95449993133545, 41.94715559232957, 3.568726936736319, 0.0821293428115375),
    @(82.84075012207032, 76.11509990692139, 52.60227175983852, 3.626715821975181, 0.1116664211506126),
    @(83.10575027465821, 76.30269992828369, 57.03749658271517, 3.864082544896479, 0.1589229220701076),
    @(83.36525020599365, 76.51909992218017, 55.61904685014895, 4.107682368461616, 0.09223540229728933),
    @(83.57775020599365, 76.72309993743896, 51.80517312587544, 4.183193892801236, 0.09491349645901925),
    @(83.59100017547607, 76.82839992523193, 40.8894915209647, 3.992264812935886, -0.04811161930917729),
    @(83.66150016784668, 77.03139995574951, 44.84915308212367, 3.71485632323183, -0.02212488395006229),
    @(83.76750011444092, 77.23459995269775, 47.54289676063065, 3.650751689194498, -0.0644647319519206),
    @(83.84950008392335, 77.44609992980958, 51.95256264927125, 3.667516880191039, -0.05852386668854559),
    @(83.85550003051758, 77.70339992523193, 53.80055688312157, 3.625168273640944, -0.04645531734059793),
    @(83.79325008392334, 77.91949993133545, 52.14447479190592, 3.557500373388179, 0.01996464295223821),
    @(83.72825012207031, 78.14899993896485, 51.30507497649368, 3.576998871603226, 0.01209921332774422),
    @(83.70825004577637, 78.35069992065429, 49.49100980073005, 3.556276933230153, 0.03165325655515838),
    @(83.66500015258789, 78.54599990844727, 56.28938443481556, 3.651115444205833, 0.03414697643152742)
)

$startRow = 183
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($startRow + $i, 1 + $j).Value = $row[$j]
    }
}

